$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue "D2" "248.43"
Set-TextValue "D3" "21.73"
Set-TextValue "D4" "5.501"
Set-TextValue "D5" "0.05637"
Set-TextValue "D7" "6.450"
Set-TextValue "D8" "0.8023"
Set-TextValue "D10" "0.1427"
Set-TextValue "D11" "0.07236"
Set-TextValue "D12" "0.03169"
Set-TextValue "D14" "0.09272"
Set-TextValue "D15" "0.001667"
Set-TextValue "D16" "3.286"
Set-TextValue "D17" "0.04750"
Set-TextValue "D18" "0.0005825"
Set-TextValue "E18" "17OneONEWorstin24h"
Set-TextValue "D19" "0.006444"
Set-TextValue "D20" "0.005024"
Set-TextValue "D22" "0.0001505"
Set-TextValue "D24" "4.075"
Set-TextValue "D25" "2.110"
Set-TextValue "D40" "0.04097"
Set-TextValue "B41" "KickToken"
Set-TextValue "C41" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006953"
Set-TextValue "E41" "40KickTokenKICK"
Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1040"
Set-TextValue "E42" "41BKEXTokenBKK"
Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003208"
Set-TextValue "E43" "42CEJICEJI"
Set-TextValue "D44" "0.009019"
Set-TextValue "D45" "0.00005656"
Set-TextValue "D47" "0.7872"
Set-TextValue "D48" "0.01704"
Set-TextValue "D49" "0.00002107"
Set-TextValue "D50" "0.01013"
